$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 54
$ws.Cells.Item($row, 1).Value = "Golang Developer"
$ws.Cells.Item($row, 2).Value = "https://www.dice.com/job-detail/5db1be0f-a4aa-4747-bd97-f8b8bce91482"
$ws.Cells.Item($row, 3).Value = "Plano, Texas"
$ws.Cells.Item($row, 4).Value = "Contract"
$ws.Cells.Item($row, 5).Value = "$63.7 - $73.76"
$ws.Cells.Item($row, 6).Value = "Judge Group, Inc."
